# Weekly update: a new price-report row for "Feria Lagunitas de Puerto Montt - Perejil"
# was added to the source data. In the canonical sheet this lands as a brand-new
# row 443 (dated 2023-12-05 / serial 45265), pushing every existing row at/after
# the old row 443 down by one (old 443 -> 444, ..., old 474 -> 475).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 443 - Excel shifts rows 443:474 down to 444:475,
# preserving all of their existing values/formatting automatically.
$ws.Rows.Item(443).Insert()

# Populate the newly-inserted row 443 with the new record.
$ws.Range("A443").Value = 4
$ws.Range("B443").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C443").Value = "Los Lagos"
$ws.Range("D443").Value = 45265
$ws.Range("E443").Value = 10
$ws.Range("F443").Value = 100112044
$ws.Range("G443").Value = "Perejil"
$ws.Range("H443").Value = "Sin especificar"
$ws.Range("I443").Value = "Primera"
$ws.Range("J443").Value = 180
$ws.Range("K443").Value = 8000
$ws.Range("L443").Value = 8000
$ws.Range("M443").Value = 8000
$ws.Range("N443").Value = "$/docena de atados (3 kilos)"
$ws.Range("O443").Value = "Región Metropolitana"
$ws.Range("P443").Value = 2667
$ws.Range("Q443").Value = 3
$ws.Range("R443").Value = "Hortaliza"
